$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Shift existing blog references down and insert the new "ser: 119" entry
$ws.Range("I8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 117"
$ws.Range("E8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 118"
$ws.Range("C8").Value = "type: blog`nwidth: 2`nheight: 1`nser: 119"
